# Update "想去人数" (interest-count) figures in column F across all four
# sheets to the freshly scraped values (gh-pages data refresh @ 456a3b4).
$wb = $excel.ActiveWorkbook

# --- 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 37514  # was 37432
$ws.Range("F4").Value = 635  # was 634
$ws.Range("F7").Value = 361  # was 360
$ws.Range("F9").Value = 839  # was 837
$ws.Range("F10").Value = 91  # was 90
$ws.Range("F11").Value = 706  # was 701
$ws.Range("F12").Value = 533  # was 532
$ws.Range("F13").Value = 33  # was 23
$ws.Range("F15").Value = 13  # was 9
$ws.Range("F16").Value = 644  # was 638
$ws.Range("F19").Value = 441  # was 439
$ws.Range("F20").Value = 1163  # was 1159
$ws.Range("F21").Value = 92  # was 91
$ws.Range("F22").Value = 820  # was 815
$ws.Range("F23").Value = 2510  # was 2505
$ws.Range("F24").Value = 995  # was 992
$ws.Range("F25").Value = 560  # was 556
$ws.Range("F27").Value = 1155  # was 1153
$ws.Range("F29").Value = 763  # was 758
$ws.Range("F30").Value = 54  # was 51
$ws.Range("F31").Value = 1149  # was 1147

# --- 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 382  # was 375
$ws.Range("F5").Value = 1  # was 0
$ws.Range("F7").Value = 56  # was 55
$ws.Range("F10").Value = 10  # was 9
$ws.Range("F11").Value = 7  # was 6

# --- 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 624  # was 621

# --- 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 624  # was 621
$ws.Range("F3").Value = 37514  # was 37432
$ws.Range("F5").Value = 635  # was 634
$ws.Range("F9").Value = 361  # was 360
$ws.Range("F11").Value = 382  # was 375
$ws.Range("F13").Value = 1  # was 0
$ws.Range("F15").Value = 839  # was 837
$ws.Range("F16").Value = 91  # was 90
$ws.Range("F17").Value = 706  # was 701
$ws.Range("F18").Value = 533  # was 532
$ws.Range("F19").Value = 33  # was 23
$ws.Range("F20").Value = 56  # was 55
$ws.Range("F24").Value = 10  # was 9
$ws.Range("F25").Value = 13  # was 9
$ws.Range("F26").Value = 7  # was 6
$ws.Range("F27").Value = 644  # was 639
$ws.Range("F30").Value = 441  # was 439
$ws.Range("F31").Value = 1163  # was 1159
$ws.Range("F32").Value = 92  # was 91
$ws.Range("F33").Value = 820  # was 815
$ws.Range("F34").Value = 2510  # was 2505
$ws.Range("F35").Value = 995  # was 992
$ws.Range("F36").Value = 560  # was 556
$ws.Range("F38").Value = 1155  # was 1153
$ws.Range("F41").Value = 763  # was 758
$ws.Range("F42").Value = 54  # was 51
$ws.Range("F43").Value = 1149  # was 1147
